$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Some new price strings look like plain numbers (e.g. "587.22"); force the cell
# to Text format first so Excel keeps them as literal strings (matches source data,
# which stores every price as text, including thousand-dot-grouped values like
# "66.898.56" that are NOT valid numbers anyway).
$prices = [ordered]@{
    "D2" = '66.898.56'
    "D3" = '3.440.26'
    "D5" = '587.22'
    "D6" = '181.22'
    "D7" = '0.630'
    "D9" = '3.436.55'
    "D13" = '4.038.10'
    "D15" = '30.04'
    "D16" = '66.830.33'
    "D18" = '3.468.33'
    "D19" = '5.96'
    "D21" = '373.08'
    "D22" = '7.67'
    "D23" = '73.38'
    "D24" = '0.0000130'
    "D26" = '0.538'
    "D27" = '10.00'
    "D29" = '1.00'
    "D30" = '5.90'
    "D32" = '23.69'
    "D33" = '0.999'
    "D34" = '7.13'
    "D37" = '162.90'
    "D38" = '0.881'
    "D39" = '27.86'
    "D41" = '2.67'
    "D42" = '4.50'
    "D45" = '0.0699'
    "D46" = '25.75'
    "D47" = '339.22'
    "D48" = '40.20'
    "D51" = '32.08'
    "D43" = '6.47'
    "D44" = '2.742.91'
}
foreach ($cell in $prices.Keys) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $prices[$cell]
}

# --- Volume(1h) % (column E) updates ---
$volumes = [ordered]@{
    "E2" = '  -0.54%  '
    "E3" = '  -1.44%  '
    "E4" = '  +0.01%  '
    "E5" = '  -0.15%  '
    "E6" = '  +1.90%  '
    "E7" = '  +4.66%  '
    "E8" = '  -0.02%  '
    "E9" = '  -1.38%  '
    "E10" = '  -0.06%  '
    "E11" = '  +0.64%  '
    "E12" = '  -1.28%  '
    "E13" = '  -1.31%  '
    "E14" = '  +1.33%  '
    "E15" = '  -2.15%  '
    "E16" = '  -0.52%  '
    "E17" = '  +0.72%  '
    "E18" = '  -0.65%  '
    "E19" = '  -1.01%  '
    "E20" = '  -0.38%  '
    "E21" = '  -2.79%  '
    "E22" = '  -2.70%  '
    "E23" = '  +0.65%  '
    "E24" = '  +6.74%  '
    "E25" = '  -0.19%  '
    "E26" = '  -1.38%  '
    "E27" = '  +0.99%  '
    "E28" = '  +1.99%  '
    "E29" = '  +0.06%  '
    "E30" = '  -0.62%  '
    "E31" = '  -0.38%  '
    "E32" = '  -3.13%  '
    "E33" = '  -0.04%  '
    "E34" = '  -1.29%  '
    "E35" = '  -3.78%  '
    "E36" = '  -0.83%  '
    "E37" = '  +1.40%  '
    "E38" = '  -1.83%  '
    "E39" = '  -6.46%  '
    "E40" = '  +0.49%  '
    "E41" = '  +0.67%  '
    "E42" = '  -0.82%  '
    "E45" = '  -0.39%  '
    "E46" = '  +3.74%  '
    "E47" = '  +7.06%  '
    "E48" = '  -1.29%  '
    "E49" = '  -3.13%  '
    "E50" = '  +2.74%  '
    "E51" = '  +2.39%  '
    "E43" = '  +0.08%  '
    "E44" = '  -0.50%  '
}
foreach ($cell in $volumes.Keys) {
    $ws.Range($cell).Value = $volumes[$cell]
}

# --- Rows 43/44 swap: RenderToken now ranks above Maker ---
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
